# Mark the first batch of tasks (rows 3-8, "Fait?" column E) as done ("oui"),
# matching the commit "4-5 taches presque terminees".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "oui"
$ws.Range("E4").Value = "oui"
$ws.Range("E5").Value = "oui"
$ws.Range("E6").Value = "oui"
$ws.Range("E7").Value = "oui"
$ws.Range("E8").Value = "oui"

# Reflect the author's updated scroll position / active cell selection.
$ws.Range("F10").Select() | Out-Null
